$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.999.55'
$ws.Range('E2').Value = '  +1.77%  '
$ws.Range('D3').Value = '1.905.85'
$ws.Range('E3').Value = '  +2.00%  '
$rng = $ws.Range('D4')
$rng.NumberFormat = "@"
$rng.Value = '1.006'
$ws.Range('E4').Value = '  -0.76%  '
$rng = $ws.Range('D5')
$rng.NumberFormat = "@"
$rng.Value = '315.95'
$ws.Range('E5').Value = '  +1.29%  '
$rng = $ws.Range('D6')
$rng.NumberFormat = "@"
$rng.Value = '1.005'
$ws.Range('E6').Value = '  -0.79%  '
$rng = $ws.Range('D7')
$rng.NumberFormat = "@"
$rng.Value = '0.4815'
$ws.Range('E7').Value = '  +0.74%  '
$rng = $ws.Range('D8')
$rng.NumberFormat = "@"
$rng.Value = '0.3804'
$ws.Range('E8').Value = '  +1.84%  '
$rng = $ws.Range('D9')
$rng.NumberFormat = "@"
$rng.Value = '0.07350'
$ws.Range('E9').Value = '  +0.51%  '
$rng = $ws.Range('D10')
$rng.NumberFormat = "@"
$rng.Value = '0.9320'
$ws.Range('E10').Value = '  -0.41%  '
$rng = $ws.Range('D11')
$rng.NumberFormat = "@"
$rng.Value = '20.78'
$ws.Range('E11').Value = '  +0.48%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.889.74'
$ws.Range('E12').Value = '  +1.25%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$rng = $ws.Range('D13')
$rng.NumberFormat = "@"
$rng.Value = '0.07749'
$ws.Range('E13').Value = '  -1.09%  '
$rng = $ws.Range('D14')
$rng.NumberFormat = "@"
$rng.Value = '5.498'
$ws.Range('E14').Value = '  +1.04%  '
$rng = $ws.Range('D15')
$rng.NumberFormat = "@"
$rng.Value = '6.659'
$ws.Range('E15').Value = '  +1.60%  '
$rng = $ws.Range('D16')
$rng.NumberFormat = "@"
$rng.Value = '91.74'
$ws.Range('E16').Value = '  +1.77%  '
$rng = $ws.Range('D17')
$rng.NumberFormat = "@"
$rng.Value = '1.006'
$ws.Range('E17').Value = '  -0.74%  '
$rng = $ws.Range('D18')
$rng.NumberFormat = "@"
$rng.Value = '0.000008821'
$ws.Range('E18').Value = '  -0.92%  '
$ws.Range('E19').Value = '  -0.73%  '
$ws.Range('D20').Value = '28.038.33'
$rng = $ws.Range('D21')
$rng.NumberFormat = "@"
$rng.Value = '14.78'
$ws.Range('E21').Value = '  +1.00%  '
$rng = $ws.Range('D22')
$rng.NumberFormat = "@"
$rng.Value = '5.170'
$ws.Range('E22').Value = '  +1.01%  '
$ws.Range('D23').Value = '2.155.29'
$ws.Range('E23').Value = '  +2.66%  '
$ws.Range('E24').Value = '  +1.69%  '
$rng = $ws.Range('D25')
$rng.NumberFormat = "@"
$rng.Value = '155.82'
$ws.Range('E25').Value = '  +0.91%  '
$rng = $ws.Range('D26')
$rng.NumberFormat = "@"
$rng.Value = '1.917'
$ws.Range('E26').Value = '  -1.65%  '
$ws.Range('E27').Value = '  +0.15%  '
$rng = $ws.Range('D28')
$rng.NumberFormat = "@"
$rng.Value = '2.123'
$ws.Range('E28').Value = '  +5.02%  '
$rng = $ws.Range('D29')
$rng.NumberFormat = "@"
$rng.Value = '116.78'
$ws.Range('E29').Value = '  +0.97%  '
$ws.Range('E30').Value = '  -0.57%  '
$rng = $ws.Range('D31')
$rng.NumberFormat = "@"
$rng.Value = '0.08934'
$ws.Range('E31').Value = '  +0.25%  '
$rng = $ws.Range('D32')
$rng.NumberFormat = "@"
$rng.Value = '3.304'
$ws.Range('E32').Value = '  -0.97%  '
$rng = $ws.Range('D34')
$rng.NumberFormat = "@"
$rng.Value = '0.7729'
$ws.Range('E34').Value = '  +1.92%  '
$ws.Range('E35').Value = '  +1.28%  '
$rng = $ws.Range('D36')
$rng.NumberFormat = "@"
$rng.Value = '2.612'
$ws.Range('E36').Value = '  -4.51%  '
$rng = $ws.Range('D37')
$rng.NumberFormat = "@"
$rng.Value = '0.02057'
$ws.Range('E37').Value = '  +0.73%  '
$rng = $ws.Range('D38')
$rng.NumberFormat = "@"
$rng.Value = '1.113'
$ws.Range('E38').Value = '  -0.53%  '
$rng = $ws.Range('D39')
$rng.NumberFormat = "@"
$rng.Value = '0.05293'
$ws.Range('E39').Value = '  +0.45%  '
$rng = $ws.Range('D40')
$rng.NumberFormat = "@"
$rng.Value = '0.5481'
$ws.Range('E40').Value = '  +3.22%  '
$ws.Range('E41').Value = '  -0.21%  '
$rng = $ws.Range('D42')
$rng.NumberFormat = "@"
$rng.Value = '7.031'
$ws.Range('E42').Value = '  -0.56%  '
$rng = $ws.Range('D43')
$rng.NumberFormat = "@"
$rng.Value = '0.1531'
$ws.Range('E43').Value = '  +0.52%  '
$rng = $ws.Range('D44')
$rng.NumberFormat = "@"
$rng.Value = '8.485'
$ws.Range('E44').Value = '  +0.08%  '
$rng = $ws.Range('D45')
$rng.NumberFormat = "@"
$rng.Value = '10.69'
$ws.Range('E45').Value = '  +0.81%  '
$rng = $ws.Range('D46')
$rng.NumberFormat = "@"
$rng.Value = '0.4823'
$ws.Range('E46').Value = '  +0.45%  '
$rng = $ws.Range('D47')
$rng.NumberFormat = "@"
$rng.Value = '107.84'
$ws.Range('E47').Value = '  +4.95%  '
$rng = $ws.Range('D48')
$rng.NumberFormat = "@"
$rng.Value = '1.005'
$ws.Range('E48').Value = '  -0.87%  '
$rng = $ws.Range('D49')
$rng.NumberFormat = "@"
$rng.Value = '1.646'
$ws.Range('E49').Value = '  -0.53%  '
$rng = $ws.Range('D50')
$rng.NumberFormat = "@"
$rng.Value = '67.88'
$ws.Range('E50').Value = '  +0.81%  '
$rng = $ws.Range('D51')
$rng.NumberFormat = "@"
$rng.Value = '0.06074'
$ws.Range('E51').Value = '  -0.14%  '
